$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = '25.181.98'
$ws.Range("E2").Value = '  -2.64%  '
$ws.Range("D3").Value = '1.660.93'
$ws.Range("E3").Value = '  -4.52%  '
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  +0.75%  '
$ws.Range("D5").Value = '237.43'
$ws.Range("E5").Value = '  -1.13%  '
$ws.Range("D6").Value = '1.008'
$ws.Range("E6").Value = '  +0.82%  '
$ws.Range("D7").Value = '0.4727'
$ws.Range("E7").Value = '  -9.54%  '
$ws.Range("D8").Value = '0.2615'
$ws.Range("E8").Value = '  -4.77%  '
$ws.Range("D9").Value = '0.05967'
$ws.Range("E9").Value = '  -3.28%  '
$ws.Range("D10").Value = '0.07114'
$ws.Range("E10").Value = '  -1.02%  '
$ws.Range("D11").Value = '1.695.40'
$ws.Range("E11").Value = '  -2.58%  '
$ws.Range("D12").Value = '0.6203'
$ws.Range("E12").Value = '  -3.65%  '
$ws.Range("D13").Value = '14.39'
$ws.Range("E13").Value = '  -3.99%  '
$ws.Range("D14").Value = '4.611'
$ws.Range("E14").Value = '  -0.29%  '
$ws.Range("D15").Value = '72.89'
$ws.Range("E15").Value = '  -5.95%  '
$ws.Range("D16").Value = '1.010'
$ws.Range("E16").Value = '  +0.96%  '
$ws.Range("D17").Value = '1.008'
$ws.Range("E17").Value = '  +0.86%  '
$ws.Range("D18").Value = '25.187.57'
$ws.Range("E18").Value = '  -2.75%  '
$ws.Range("D19").Value = '11.40'
$ws.Range("E19").Value = '  -2.79%  '
$ws.Range("D20").Value = '0.000006562'
$ws.Range("E20").Value = '  -3.23%  '
$ws.Range("D21").Value = '1.902.93'
$ws.Range("E21").Value = '  -3.28%  '
$ws.Range("D22").Value = '4.424'
$ws.Range("E22").Value = '  +3.14%  '
$ws.Range("D23").Value = '8.552'
$ws.Range("E23").Value = '  -1.13%  '
$ws.Range("D24").Value = '5.250'
$ws.Range("E24").Value = '  -0.56%  '
$ws.Range("D25").Value = '133.27'
$ws.Range("E25").Value = '  -3.90%  '
$ws.Range("E26").Value = '  -3.27%  '
$ws.Range("D27").Value = '1.359'
$ws.Range("E27").Value = '  -11.11%  '
$ws.Range("D28").Value = '1.705'
$ws.Range("E28").Value = '  -3.58%  '
$ws.Range("D29").Value = '102.40'
$ws.Range("E29").Value = '  -3.26%  '
$ws.Range("D30").Value = '3.829'
$ws.Range("E30").Value = '  -2.47%  '
$ws.Range("D31").Value = '0.07874'
$ws.Range("E31").Value = '  -5.37%  '
$ws.Range("D32").Value = '3.530'
$ws.Range("E32").Value = '  -4.32%  '
$ws.Range("D33").Value = '0.04605'
$ws.Range("E33").Value = '  -0.63%  '
$ws.Range("D34").Value = '2.633'
$ws.Range("E34").Value = '  -0.26%  '
$ws.Range("D35").Value = '0.9392'
$ws.Range("E35").Value = '  -5.13%  '
$ws.Range("D36").Value = '0.5795'
$ws.Range("E36").Value = '  -6.46%  '
$ws.Range("D37").Value = '2.609'
$ws.Range("E37").Value = '  -2.78%  '
$ws.Range("D38").Value = '0.01545'
$ws.Range("E38").Value = '  -3.83%  '
$ws.Range("D39").Value = '1.008'
$ws.Range("E39").Value = '  +0.94%  '
$ws.Range("D40").Value = '0.8382'
$ws.Range("E40").Value = '  +13.30%  '
$ws.Range("D41").Value = '1.831'
$ws.Range("E41").Value = '  -5.68%  '
$ws.Range("D42").Value = '98.72'
$ws.Range("E42").Value = '  +0.85%  '
$ws.Range("D43").Value = '0.3702'
$ws.Range("E43").Value = '  -4.01%  '
$ws.Range("D44").Value = '4.867'
$ws.Range("E44").Value = '  -2.57%  '
$ws.Range("D45").Value = '0.1134'
$ws.Range("E45").Value = '  +0.19%  '
$ws.Range("D46").Value = '0.05168'
$ws.Range("E46").Value = '  -1.47%  '
$ws.Range("D47").Value = '6.032'
$ws.Range("E47").Value = '  -3.58%  '
$ws.Range("D48").Value = '53.43'
$ws.Range("E48").Value = '  -2.40%  '
$ws.Range("D49").Value = '29.64'
$ws.Range("E49").Value = '  -2.77%  '
$ws.Range("D50").Value = '1.007'
$ws.Range("E50").Value = '  +0.63%  '
$ws.Range("D51").Value = '7.362'
$ws.Range("E51").Value = '  -3.33%  '

$rng.ClearFormats()
